# Apply the reshuffle of rows 3-7 described by the commit diff.
# The "Id" (A), "Ost" (Q) and "Nord" (R) columns are rotated among rows 3-7,
# while the taxon-describing columns (B, D, E, F, G, H, M) are swapped
# between rows 3 and 6. Rather than trying to replicate the move
# operation cell-by-cell (which risks clobbering values before they are
# read), we just write the final, known target values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (Id), Q (Ost), R (Nord) - final values per row ---
$ws.Range("A3").Value = 111742281
$ws.Range("Q3").Value = 331821.5503750234
$ws.Range("R3").Value = 6626517.909892835

$ws.Range("A4").Value = 111742294
$ws.Range("Q4").Value = 331799.9927276275
$ws.Range("R4").Value = 6626510.806996167

$ws.Range("A5").Value = 111742278
$ws.Range("Q5").Value = 331818.8411813352
$ws.Range("R5").Value = 6626525.099085328

$ws.Range("A6").Value = 111742269
$ws.Range("Q6").Value = 331779.9179887357
$ws.Range("R6").Value = 6626525.342625097

$ws.Range("A7").Value = 111742299
$ws.Range("Q7").Value = 331807.7707727421
$ws.Range("R7").Value = 6626503.893626045

# --- Row 3 becomes the "Thomsons trägnagare" record ---
$ws.Range("B3").Value = 4711
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 100299
$ws.Range("F3").Value = "Thomsons trägnagare"
$ws.Range("G3").Value = "Cacotemnus thomsoni"
$ws.Range("H3").Value = "(Kraatz, 1881)"
$ws.Range("M3").Value = "färska gnagspår"

# --- Row 6 becomes a "Vedtrappmossa" record (same as rows 4,5,7) ---
$ws.Range("B6").Value = 94134
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 53
$ws.Range("F6").Value = "Vedtrappmossa"
$ws.Range("G6").Value = "Crossocalyx hellerianus"
$ws.Range("H6").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("M6").ClearContents()
